$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values
$ws.Range("B2").Value = 15.458725709495777
$ws.Range("C2").Value = 10.787023792862547
$ws.Range("D2").Value = 15.923135859821482
$ws.Range("E2").Value = 8.7728464407593894

# Row 3 (STR) values
$ws.Range("B3").Value = 13.121846357517555
$ws.Range("C3").Value = 13.299741894026944
$ws.Range("D3").Value = 12.754086142636629
$ws.Range("E3").Value = 15.748702774842679

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
